$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - headers
$ws.Range("A1").Value = "Grado, código"
$ws.Range("B1").Value = "Personas"
$ws.Range("C1").Value = "Municipio codigo"
$ws.Range("D1").Value = "Grado"
$ws.Range("E1").Value = "Sexo, código"
$ws.Range("F1").Value = "Sexo"
$ws.Range("G1").Value = "Municipio nombre"

# Row 2
$ws.Range("A2").Value = "null"
$ws.Range("B2").Value = "iaest-measure:personas"
$ws.Range("C2").Value = "null"
$ws.Range("D2").Value = "iaest-measure:grado"
$ws.Range("E2").Value = "null"
$ws.Range("F2").Value = "iaest-measure:sexo"
$ws.Range("G2").Value = "sdmx-dimension:refArea"

# Row 3
$ws.Range("A3").Value = "null"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "null"
$ws.Range("D3").Value = "medida"
$ws.Range("E3").Value = "null"
$ws.Range("F3").Value = "medida"
$ws.Range("G3").Value = "dim"

# Row 4
$ws.Range("A4").Value = "null"
$ws.Range("B4").Value = "xsd:double"
$ws.Range("C4").Value = "null"
$ws.Range("D4").Value = "xsd:string"
$ws.Range("E4").Value = "null"
$ws.Range("F4").Value = "xsd:string"
$ws.Range("G4").Value = "URI-Municipio"
